$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing month-4 (April/2025) rows with revised totals
$ws.Range("C248").Value = 271101.06
$ws.Range("C249").Value = 97761
$ws.Range("C250").Value = 116177.05
$ws.Range("C251").Value = 82407.5

# Add new month-5 (May/2025) rows for the four stores
$ws.Range("A252").Value = 5
$ws.Range("B252").Value = 1
$ws.Range("C252").Value = 52798.18
$ws.Range("D252").Value = 2025
$ws.Range("E252").Value = "Bibi Cell Mundi"

$ws.Range("A253").Value = 5
$ws.Range("B253").Value = 2
$ws.Range("C253").Value = 18395.5
$ws.Range("D253").Value = 2025
$ws.Range("E253").Value = "Bibi Cell Manauara"

$ws.Range("A254").Value = 5
$ws.Range("B254").Value = 3
$ws.Range("C254").Value = 39550.8
$ws.Range("D254").Value = 2025
$ws.Range("E254").Value = "Bibi Cell Vieiralves"

$ws.Range("A255").Value = 5
$ws.Range("B255").Value = 4
$ws.Range("C255").Value = 15862.08
$ws.Range("D255").Value = 2025
$ws.Range("E255").Value = "Bibi Cell Ponta Negra"
